$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL, Version, Date, Publisher values ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/path"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": move the ele-1/ext-1 constraint text from the
#     "Extension" row (row 2) to the "Extension.extension" row (row 4),
#     in the "Constraint(s)" column (AI) ---
$wsElem = $wb.Worksheets.Item("Elements")
$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$wsElem.Range("AI2").Value = ""
$wsElem.Range("AI4").Value = $constraintText
